# Consolidate the title's separate text runs ("Two-Column", " ", "Layout")
# into a single run containing "Two-Column Layout".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# Assigning the TextRange.Text directly to the already-equivalent
# concatenated string is treated as a no-op by the writer (the rendered
# text doesn't change) and leaves the separate runs untouched. It also does
# a prefix-preserving diff, so a temporary value that shares a prefix with
# the target text would only grow/trim the existing runs instead of
# consolidating them. Using a placeholder with no shared prefix forces a
# real text replacement, collapsing the paragraph down to a single run;
# then we set the real desired text onto that single run.
$tr.Text = "Z"
$tr.Text = "Two-Column Layout"
